# 54BND_svm.xlsx update:
#  - drop the "0,001" sheet, keep/rename the remaining sheet to "54BND"
#  - add a third kernel-group ("0,05") of results in columns H:J
#  - recompute the accuracy row, rename std_dev -> sensibility, add specificity row

$wb = $excel.ActiveWorkbook

# ---- 1. Remove the "0,001" sheet and rename the remaining one ----
$wb.Worksheets.Item("0,001").Delete() | Out-Null
$ws = $wb.ActiveSheet
$ws.Name = "54BND"

# ---- 2. Row 1 headers: add rbf/poly/linear for the new 0,05 group ----
$ws.Range("B1:D1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$ws.Range("H1").Value = "rbf"
$ws.Range("I1").Value = "poly"
$ws.Range("J1").Value = "linear"

# ---- 3. Row 2 group label: merge H2:J2 (before formatting, so the merge
#         does not recompute "smart" partial borders) then copy the format
#         from the existing E2:G2 group and set its text ----
$ws.Range("H2:J2").MergeCells = $true
$ws.Range("E2:G2").Copy()
$ws.Range("H2:J2").PasteSpecial(-4122)
$ws.Range("H2").Value = "0,05"

# ---- 4. Fold results (rows 3-22) for the new 0,05 kernel group ----
$ws.Range("H3").Value = 0.6666666666666666; $ws.Range("I3").Value = 0.6666666666666666; $ws.Range("J3").Value = 0
$ws.Range("H4").Value = 0.6666666666666666; $ws.Range("I4").Value = 0.6666666666666666; $ws.Range("J4").Value = 0
$ws.Range("H5").Value = 0; $ws.Range("I5").Value = 0; $ws.Range("J5").Value = 0
$ws.Range("H6").Value = 0.3333333333333333; $ws.Range("I6").Value = 0.3333333333333333; $ws.Range("J6").Value = 0.3333333333333333
$ws.Range("H7").Value = 0.6666666666666666; $ws.Range("I7").Value = 0.6666666666666666; $ws.Range("J7").Value = 0.6666666666666666
$ws.Range("H8").Value = 0.6666666666666666; $ws.Range("I8").Value = 0.6666666666666666; $ws.Range("J8").Value = 0.6666666666666666
$ws.Range("H9").Value = 0.6666666666666666; $ws.Range("I9").Value = 0.6666666666666666; $ws.Range("J9").Value = 0.6666666666666666
$ws.Range("H10").Value = 0.6666666666666666; $ws.Range("I10").Value = 0.6666666666666666; $ws.Range("J10").Value = 0.6666666666666666
$ws.Range("H11").Value = 0.6666666666666666; $ws.Range("I11").Value = 0.6666666666666666; $ws.Range("J11").Value = 0.3333333333333333
$ws.Range("H12").Value = 0.6666666666666666; $ws.Range("I12").Value = 0.6666666666666666; $ws.Range("J12").Value = 0.6666666666666666
$ws.Range("H13").Value = 0.3333333333333333; $ws.Range("I13").Value = 0.3333333333333333; $ws.Range("J13").Value = 0.3333333333333333
$ws.Range("H14").Value = 0.3333333333333333; $ws.Range("I14").Value = 0.3333333333333333; $ws.Range("J14").Value = 0.6666666666666666
$ws.Range("H15").Value = 0; $ws.Range("I15").Value = 0; $ws.Range("J15").Value = 0.3333333333333333
$ws.Range("H16").Value = 0.3333333333333333; $ws.Range("I16").Value = 0.3333333333333333; $ws.Range("J16").Value = 0.3333333333333333
$ws.Range("H17").Value = 0; $ws.Range("I17").Value = 0.5; $ws.Range("J17").Value = 0
$ws.Range("H18").Value = 0; $ws.Range("I18").Value = 0; $ws.Range("J18").Value = 0
$ws.Range("H19").Value = 0.5; $ws.Range("I19").Value = 0.5; $ws.Range("J19").Value = 0
$ws.Range("H20").Value = 1; $ws.Range("I20").Value = 1; $ws.Range("J20").Value = 0.5
$ws.Range("H21").Value = 0; $ws.Range("I21").Value = 0; $ws.Range("J21").Value = 0
$ws.Range("H22").Value = 0; $ws.Range("I22").Value = 0; $ws.Range("J22").Value = 0

# ---- 5. Row 23 ("accuracy"): recomputed values for every kernel group ----
$ws.Range("B23").Value = 0.7407407407407407; $ws.Range("C23").Value = 0.6666666666666666; $ws.Range("D23").Value = 0.7222222222222222; $ws.Range("E23").Value = 0.7592592592592593; $ws.Range("F23").Value = 0.6851851851851852; $ws.Range("G23").Value = 0.7407407407407407; $ws.Range("H23").Value = 0.4259259259259259; $ws.Range("I23").Value = 0.4444444444444444; $ws.Range("J23").Value = 0.3333333333333333

# ---- 6. Row 24: relabel std_dev -> sensibility, with its new values ----
$ws.Range("A24").Value = "sensibility"
$ws.Range("B24").Value = 0.5714285714285714; $ws.Range("C24").Value = 0.7857142857142857; $ws.Range("D24").Value = 0.5357142857142857; $ws.Range("E24").Value = 0.6071428571428571; $ws.Range("F24").Value = 0.75; $ws.Range("G24").Value = 0.6071428571428571; $ws.Range("H24").Value = 0.8214285714285714; $ws.Range("I24").Value = 0.8571428571428571; $ws.Range("J24").Value = 0.6071428571428571

# ---- 7. Row 25 (new): "specificity" row, formatted like row 24's label ----
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = "specificity"
$ws.Range("B25").Value = 0.9230769230769231; $ws.Range("C25").Value = 0.5384615384615384; $ws.Range("D25").Value = 0.9230769230769231; $ws.Range("E25").Value = 0.9230769230769231; $ws.Range("F25").Value = 0.6153846153846154; $ws.Range("G25").Value = 0.8846153846153846; $ws.Range("H25").Value = 0; $ws.Range("I25").Value = 0; $ws.Range("J25").Value = 0.03846153846153846

# ---- 8. Selection matches the saved cursor position ----
$ws.Range("I12").Select() | Out-Null
